$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6707215716073449
$ws.Range("C2").Value = 0.1579930342811906
$ws.Range("D2").Value = 0.01294316665136108
$ws.Range("E2").Value = 0.1264417977370798
$ws.Range("F2").Value = 0.5538072680917523
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("M2").Value = 0.309925283911106
$ws.Range("N2").Value = 0.9722120263364395
$ws.Range("O2").Value = 1.797013954961699
$ws.Range("B3").Value = 0.5892113537390742
$ws.Range("C3").Value = 0.1408492153979637
$ws.Range("D3").Value = 0.01180183988999062
$ws.Range("E3").Value = 0.1199290195625125
$ws.Range("F3").Value = 0.5440222806960406
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("M3").Value = 0.2762130393736442
$ws.Range("N3").Value = 0.9838772478457471
$ws.Range("O3").Value = 1.778239503835522
$ws.Range("B4").Value = 0.5391197371739906
$ws.Range("C4").Value = 0.1302612413028328
$ws.Range("D4").Value = 0.01109670045016387
$ws.Range("E4").Value = 0.1160393710671954
$ws.Range("F4").Value = 0.5384551950288667
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("M4").Value = 0.2555826412098554
$ws.Range("N4").Value = 0.9914898772992515
$ws.Range("O4").Value = 1.768167780096121
$ws.Range("B5").Value = 0.5186967731641516
$ws.Range("C5").Value = 0.1259312547790046
$ws.Range("D5").Value = 0.01080827466468648
$ws.Range("E5").Value = 0.1144815140420974
$ws.Range("F5").Value = 0.5362971772272189
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("M5").Value = 0.2471929050606434
$ws.Range("N5").Value = 0.994705237091658
$ws.Range("O5").Value = 1.764428586232015
$ws.Range("B6").Value = 0.5153049630231692
$ws.Range("C6").Value = 0.1252113452438266
$ws.Range("D6").Value = 0.01076031742011807
$ws.Range("E6").Value = 0.1142244687087199
$ws.Range("F6").Value = 0.5359455132181949
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("M6").Value = 0.2458008404967487
$ws.Range("N6").Value = 0.9952459781507095
$ws.Range("O6").Value = 1.763829719546152
$ws.Range("B7").Value = 0.5388443462122723
$ws.Range("C7").Value = 0.1302029072920448
$ws.Range("D7").Value = 0.01109281496765391
$ws.Range("E7").Value = 0.1160182514466328
$ws.Range("F7").Value = 0.5384256437197976
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("M7").Value = 0.2554694242300357
$ws.Range("N7").Value = 0.9915327826713778
$ws.Range("O7").Value = 1.768115874978093
$ws.Range("B8").Value = 0.6426264640834063
$ws.Range("C8").Value = 0.1520947386414377
$ws.Range("D8").Value = 0.01255055418335616
$ws.Range("E8").Value = 0.1241733643581853
$ws.Range("F8").Value = 0.5503417257165495
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("M8").Value = 0.2982869307483895
$ws.Range("N8").Value = 0.9761407452748259
$ws.Range("O8").Value = 1.79023776193705
$ws.Range("B9").Value = 0.8457665818863234
$ws.Range("C9").Value = 0.1945294744799924
$ws.Range("D9").Value = 0.01537380291156154
$ws.Range("E9").Value = 0.1410447569340647
$ws.Range("F9").Value = 0.5772216663375076
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("M9").Value = 0.3828063636020644
$ws.Range("N9").Value = 0.9495302693902303
$ws.Range("O9").Value = 1.845219442191137
$ws.Range("B10").Value = 0.9947626019347808
$ws.Range("C10").Value = 0.2253989016716957
$ws.Range("D10").Value = 0.01742562688303906
$ws.Range("E10").Value = 0.15399495087199
$ws.Range("F10").Value = 0.5991339525288168
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("M10").Value = 0.4452561602039964
$ws.Range("N10").Value = 0.932159991699514
$ws.Range("O10").Value = 1.892761987288452
$ws.Range("B11").Value = 1.062486596629753
$ws.Range("C11").Value = 0.2393744685263073
$ws.Range("D11").Value = 0.0183540190977638
$ws.Range("E11").Value = 0.1600109248143298
$ws.Range("F11").Value = 0.6095771854180612
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("M11").Value = 0.4737466196331468
$ws.Range("N11").Value = 0.9247317984628864
$ws.Range("O11").Value = 1.915959008154346
$ws.Range("B12").Value = 1.088123353767685
$ws.Range("C12").Value = 0.2446568546791639
$ws.Range("D12").Value = 0.01870484080906465
$ws.Range("E12").Value = 0.1623072798199132
$ws.Range("F12").Value = 0.6136004393797947
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("M12").Value = 0.4845471258391285
$ws.Range("N12").Value = 0.9219871133193962
$ws.Range("O12").Value = 1.924970030011679
$ws.Range("B13").Value = 1.082602427867755
$ws.Range("C13").Value = 0.2435196400684845
$ws.Range("D13").Value = 0.01862931835413661
$ws.Range("E13").Value = 0.1618119031397782
$ws.Range("F13").Value = 0.6127309027530998
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("M13").Value = 0.4822205187690543
$ws.Range("N13").Value = 0.9225751950130245
$ws.Range("O13").Value = 1.92301924019452
$ws.Range("B14").Value = 1.064595930777159
$ws.Range("C14").Value = 0.2398092523784214
$ws.Range("D14").Value = 0.01838289637426982
$ws.Range("E14").Value = 0.1601994802296574
$ws.Range("F14").Value = 0.6099068038609232
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("M14").Value = 0.4746349473382168
$ws.Range("N14").Value = 0.9245046239218837
$ws.Range("O14").Value = 1.916695798392567
$ws.Range("B15").Value = 1.053565251177645
$ws.Range("C15").Value = 0.2375352427937401
$ws.Range("D15").Value = 0.01823185873078614
$ws.Range("E15").Value = 0.1592142078853698
$ws.Range("F15").Value = 0.6081859078227723
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("M15").Value = 0.4699901019325807
$ws.Range("N15").Value = 0.9256953409738102
$ws.Range("O15").Value = 1.912852078020649
$ws.Range("B16").Value = 0.9903354920342622
$ws.Range("C16").Value = 0.2244841974915346
$ws.Range("D16").Value = 0.01736485189061909
$ws.Range("E16").Value = 0.1536043294146765
$ws.Range("F16").Value = 0.5984610503703323
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("M16").Value = 0.4433958892783068
$ws.Range("N16").Value = 0.9326549809653102
$ws.Range("O16").Value = 1.891277687039491
$ws.Range("B17").Value = 0.9515313190291863
$ws.Range("C17").Value = 0.2164604706176476
$ws.Range("D17").Value = 0.01683167659034979
$ws.Range("E17").Value = 0.1501950311697868
$ws.Range("F17").Value = 0.5926170893340554
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("M17").Value = 0.4271021592477098
$ws.Range("N17").Value = 0.9370458749668131
$ws.Range("O17").Value = 1.878445309108457
$ws.Range("B18").Value = 0.9292070599042859
$ws.Range("C18").Value = 0.2118391327818756
$ws.Range("D18").Value = 0.0165245396403364
$ws.Range("E18").Value = 0.1482458291065711
$ws.Range("F18").Value = 0.5893005104469751
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("M18").Value = 0.4177381277398808
$ws.Range("N18").Value = 0.9396159752171513
$ws.Range("O18").Value = 1.87121214483787
$ws.Range("B19").Value = 0.9216475976370475
$ws.Range("C19").Value = 0.2102733511788983
$ws.Range("D19").Value = 0.01642046854276913
$ws.Range("E19").Value = 0.1475878699329058
$ws.Range("F19").Value = 0.5881852444154134
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("M19").Value = 0.4145689474829481
$ws.Range("N19").Value = 0.940493819613117
$ws.Range("O19").Value = 1.868788450514785
$ws.Range("B20").Value = 0.9556626267054185
$ws.Range("C20").Value = 0.2173152638370084
$ws.Range("D20").Value = 0.01688848267304621
$ws.Range("E20").Value = 0.1505567401156824
$ws.Range("F20").Value = 0.5932345592257775
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("M20").Value = 0.428835858387572
$ws.Range("N20").Value = 0.9365738428155481
$ws.Range("O20").Value = 1.879796045364685
$ws.Range("B21").Value = 1.069885123074243
$ws.Range("C21").Value = 0.2408993523891922
$ws.Range("D21").Value = 0.01845529675117774
$ws.Range("E21").Value = 0.1606725907553965
$ws.Range("F21").Value = 0.6107344452886565
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("M21").Value = 0.4768626933547466
$ws.Range("N21").Value = 0.9239360518294504
$ws.Range("O21").Value = 1.918546982150332
$ws.Range("B22").Value = 1.144484112123337
$ws.Range("C22").Value = 0.2562553624709665
$ws.Range("D22").Value = 0.01947497517278407
$ws.Range("E22").Value = 0.1673903097266063
$ws.Range("F22").Value = 0.6225717856962802
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("M22").Value = 0.5083197299137794
$ws.Range("N22").Value = 0.9160741583198089
$ws.Range("O22").Value = 1.94519551010228
$ws.Range("B23").Value = 1.104674295187863
$ws.Range("C23").Value = 0.2480649063987528
$ws.Range("D23").Value = 0.01893115670629442
$ws.Range("E23").Value = 0.1637951056822544
$ws.Range("F23").Value = 0.6162172640931658
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("M23").Value = 0.4915242176266617
$ws.Range("N23").Value = 0.9202337802595792
$ws.Range("O23").Value = 1.930851313903247
$ws.Range("B24").Value = 0.9537949098548211
$ws.Range("C24").Value = 0.2169288379442094
$ws.Range("D24").Value = 0.01686280253329642
$ws.Range("E24").Value = 0.1503931777314307
$ws.Range("F24").Value = 0.5929552665679125
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("M24").Value = 0.4280520422652501
$ws.Range("N24").Value = 0.9367871061484649
$ws.Range("O24").Value = 1.879184927966719
$ws.Range("B25").Value = 0.7908543311912695
$ws.Range("C25").Value = 0.1831033339350938
$ws.Range("D25").Value = 0.01461391480474106
$ws.Range("E25").Value = 0.1363845327678774
$ws.Range("F25").Value = 0.5695715863173945
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("M25").Value = 0.3598806770491976
$ws.Range("N25").Value = 0.9563463693884628
$ws.Range("O25").Value = 1.829095709864845
